# Apply the data edits to the "Sprint1" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# Row 2 (US01): update Actual Size, Actual Time, Completed date
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 44473

# Row 4 (US03): add Status = Done, Actual Size, Actual Time, Completed date
$ws.Range("D4").Value = "Done"
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 44472
$ws.Range("I4").NumberFormat = "d-mmm"

# Row 5 (US04): clear Status, Actual Size, Actual Time, Completed date
$ws.Range("D5").Clear()
$ws.Range("G5").Clear()
$ws.Range("H5").Clear()
$ws.Range("I5").Clear()

# Row 7 (US08): clear Status, Actual Size, Actual Time, Completed date
$ws.Range("D7").Clear()
$ws.Range("G7").Clear()
$ws.Range("H7").Clear()
$ws.Range("I7").Clear()

# Update the active cell selection on the Sprint1 sheet to match the saved view.
$ws.Range("I20").Select() | Out-Null
